# Translate the title-slide text from Greek to English ("changed to title in english").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "CustomShape 1" / shape 1: title -> Diploma Thesis
$titleShape = $s.Shapes.Item("CustomShape 1")
$titleShape.TextFrame.TextRange.Runs(1).Text = "Diploma Thesis"

# "CustomShape 7" / shape 7: author name -> Thaleia-Dimitra Doudali
$authorShape = $s.Shapes.Item("CustomShape 7")
$authorShape.TextFrame.TextRange.Runs(1).Text = "Thaleia-Dimitra Doudali"

# "CustomShape 8" / shape 8: thesis subject -> English translation
$subjectShape = $s.Shapes.Item("CustomShape 8")
$subjectShape.TextFrame.TextRange.Runs(1).Text = "Performance evaluation of social networking services using a spatio-temporal and textual Big Data generator"
